# Auto-generated Excel COM-interop script
# Applies scheduled-runner market-data refresh values per the commit diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 874.5
$ws.Range("I12").Value = 113.5
$ws.Range("K12").Value = 113.5
$ws.Range("M12").Value = 56.5
# Row 18
$ws.Range("H18").Value = 1150
$ws.Range("I18").Value = 1150
$ws.Range("K18").Value = 1150
$ws.Range("M18").Value = -866
# Row 40
$ws.Range("H40").Value = 5566.5
$ws.Range("I40").Value = 2879.4
$ws.Range("K40").Value = 2879.4
$ws.Range("M40").Value = -2704.4
# Row 92
$ws.Range("H92").Value = 1775.2858
$ws.Range("I92").Value = 768.375
$ws.Range("K92").Value = 768.375
$ws.Range("M92").Value = 479.625
# Row 112
$ws.Range("H112").Value = 2810.8333
$ws.Range("I112").Value = 587
$ws.Range("J112").Value = 3255.6
$ws.Range("K112").Value = 1761
$ws.Range("L112").Value = 9766.799999999999
$ws.Range("M112").Value = -653
$ws.Range("N112").Value = -11982.8
# Row 113
$ws.Range("H113").Value = 10670.333
$ws.Range("I113").Value = 6502.5
$ws.Range("K113").Value = 6502.5
$ws.Range("M113").Value = -3248.5
# Row 131
$ws.Range("H131").Value = 15818
$ws.Range("I131").Value = 13772.5
$ws.Range("K131").Value = 41317.5
$ws.Range("M131").Value = -36277.5
# Row 132
$ws.Range("H132").Value = 1844.0769
$ws.Range("I132").Value = 1844.0769
$ws.Range("K132").Value = 5532.2307
$ws.Range("M132").Value = -3002.2307
# Row 138
$ws.Range("H138").Value = 3928.4375
$ws.Range("J138").Value = 3738.9583
$ws.Range("L138").Value = 11216.8749
$ws.Range("N138").Value = -21496.8749

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5201.3516
$ws.Range("I32").Value = 4447.5938
$ws.Range("K32").Value = 4447.5938
$ws.Range("M32").Value = -4160.5938
# Row 45
$ws.Range("H45").Value = 47622416
$ws.Range("I45").Value = 83334760
$ws.Range("J45").Value = 5955.5557
$ws.Range("K45").Value = 83334760
$ws.Range("L45").Value = 5955.5557
$ws.Range("M45").Value = -83334383
$ws.Range("N45").Value = -6709.5557
# Row 61
$ws.Range("H61").Value = 3385.963
$ws.Range("I61").Value = 2754.6538
$ws.Range("K61").Value = 2754.6538
$ws.Range("M61").Value = -2542.6538
# Row 74
$ws.Range("H74").Value = 13890621
$ws.Range("J74").Value = 2004.3334
$ws.Range("L74").Value = 2004.3334
$ws.Range("N74").Value = -3752.3334
# Row 77
$ws.Range("H77").Value = 13890621
$ws.Range("J77").Value = 2004.3334
$ws.Range("L77").Value = 10021.667
$ws.Range("N77").Value = -18757.667
# Row 120
$ws.Range("H120").Value = 68389.75
$ws.Range("J120").Value = 68389.75
$ws.Range("L120").Value = 68389.75
$ws.Range("N120").Value = -78065.75
# Row 122
$ws.Range("H122").Value = 1730.5385
$ws.Range("I122").Value = 1399.76
$ws.Range("K122").Value = 4199.28
$ws.Range("M122").Value = -1749.28
# Row 132
$ws.Range("H132").Value = 4256.4443
$ws.Range("I132").Value = 2440.1333
$ws.Range("J132").Value = 13338
$ws.Range("K132").Value = 7320.3999
$ws.Range("L132").Value = 40014
$ws.Range("M132").Value = -4790.3999
$ws.Range("N132").Value = -45074
# Row 136
$ws.Range("H136").Value = 3385.963
$ws.Range("I136").Value = 2754.6538
$ws.Range("K136").Value = 8263.9614
$ws.Range("M136").Value = -5713.9614

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 3051.9092
$ws.Range("I86").Value = 1894.9584
$ws.Range("J86").Value = 6137.1113
$ws.Range("K86").Value = 1894.9584
$ws.Range("L86").Value = 6137.1113
$ws.Range("M86").Value = -771.9584
$ws.Range("N86").Value = -8383.1113
# Row 89
$ws.Range("H89").Value = 3051.9092
$ws.Range("I89").Value = 1894.9584
$ws.Range("J89").Value = 6137.1113
$ws.Range("K89").Value = 9474.791999999999
$ws.Range("L89").Value = 30685.5565
$ws.Range("M89").Value = -3858.791999999999
$ws.Range("N89").Value = -41917.5565
# Row 105
$ws.Range("H105").Value = 24180.715
$ws.Range("I105").Value = 38086.832
$ws.Range("K105").Value = 38086.832
$ws.Range("M105").Value = -36339.832
# Row 107
$ws.Range("H107").Value = 1351.125
$ws.Range("I107").Value = 1351.125
$ws.Range("K107").Value = 1351.125
$ws.Range("M107").Value = 568.875
# Row 124
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
# Row 125
$ws.Range("H125").Value = 49999.5
$ws.Range("J125").Value = 49999.5
$ws.Range("L125").Value = 49999.5
$ws.Range("N125").Value = -59839.5
# Row 126
$ws.Range("H126").Value = 100000
$ws.Range("J126").Value = 100000
$ws.Range("L126").Value = 100000
$ws.Range("N126").Value = -109880
# Row 127
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
# Row 134
$ws.Range("H134").Value = 3271.6667
$ws.Range("I134").Value = 1803.875
$ws.Range("K134").Value = 5411.625
$ws.Range("M134").Value = -2876.625

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 134
$ws.Range("H134").Value = 3608.238
$ws.Range("I134").Value = 2549.9412
$ws.Range("J134").Value = 8106
$ws.Range("K134").Value = 7649.823600000001
$ws.Range("L134").Value = 24318
$ws.Range("M134").Value = -5114.823600000001
$ws.Range("N134").Value = -29388

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 22
$ws.Range("H22").Value = 9001
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 9001
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 27003
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -27341
# Row 27
$ws.Range("H27").Value = 9001
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 9001
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 27003
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -27207
# Row 38
$ws.Range("H38").Value = 141.6
$ws.Range("I38").Value = 44.333332
$ws.Range("J38").Value = 287.5
$ws.Range("K38").Value = 132.999996
$ws.Range("L38").Value = 862.5
$ws.Range("M38").Value = 214.000004
$ws.Range("N38").Value = -1556.5
# Row 87
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()
# Row 90
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()
# Row 129
$ws.Range("H129").Value = 4905171.5
$ws.Range("I129").Value = 696.125
$ws.Range("J129").Value = 9264705
$ws.Range("K129").Value = 2088.375
$ws.Range("L129").Value = 27794115
$ws.Range("M129").Value = 2911.625
$ws.Range("N129").Value = -27804115

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 4749.2856
$ws.Range("I102").Value = 2372
$ws.Range("K102").Value = 2372
$ws.Range("M102").Value = -750
# Row 122
$ws.Range("H122").Value = 2364
$ws.Range("I122").Value = 1888.16
$ws.Range("K122").Value = 5664.48
$ws.Range("M122").Value = -3214.48
# Row 126
$ws.Range("H126").Value = 3684.182
$ws.Range("I126").Value = 2001.5
$ws.Range("K126").Value = 6004.5
$ws.Range("M126").Value = -3534.5
# Row 132
$ws.Range("H132").Value = 72009.92999999999
$ws.Range("I132").Value = 102923.5
$ws.Range("J132").Value = 10182.8
$ws.Range("K132").Value = 308770.5
$ws.Range("L132").Value = 30548.4
$ws.Range("M132").Value = -306240.5
$ws.Range("N132").Value = -35608.39999999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 4424.1763
$ws.Range("I40").Value = 3121.5
$ws.Range("J40").Value = 10503.333
$ws.Range("K40").Value = 3121.5
$ws.Range("L40").Value = 10503.333
$ws.Range("M40").Value = -2985.5
$ws.Range("N40").Value = -10775.333
# Row 100
$ws.Range("H100").Value = 10002
$ws.Range("I100").Value = 1000
$ws.Range("J100").Value = 19004
$ws.Range("K100").Value = 1000
$ws.Range("L100").Value = 19004
$ws.Range("M100").Value = -459
$ws.Range("N100").Value = -20086
# Row 132
$ws.Range("H132").Value = 4215.5186
$ws.Range("I132").Value = 2172.2942
$ws.Range("K132").Value = 6516.882599999999
$ws.Range("M132").Value = -3986.882599999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 3123.8572
$ws.Range("J126").Value = 7000
$ws.Range("L126").Value = 21000
$ws.Range("N126").Value = -25940
# Row 132
$ws.Range("H132").Value = 4241.4062
$ws.Range("I132").Value = 4024.8276
$ws.Range("K132").Value = 12074.4828
$ws.Range("M132").Value = -9544.4828
